# Update ID names of User Transactions review sheet:
# row 15 (Number 14, "account/user deletion criteria must be mentioned")
# Status changes from "opened" to "closed".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coaching Review")

# Update the status cell for the last row from "opened" to "closed"
$ws.Range("E15").Value = "closed"

# Move the active selection to H14 to match the saved view state
$ws.Range("H14").Select()
